# Applies the weekly refresh to the "Frambuesa" (raspberry) price sheet:
#  - rows 49-79 get updated Fecha/Calidad/Volumen/Precio/Origen values
#    (the weekly series was refreshed/re-sorted)
#  - two brand new weekly rows are inserted at 80-81
#  - the former last row (old row 80) shifts down to row 82

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: update existing rows 49-79 in place ---
$ws.Range("D49").Value2 = 44589
$ws.Range("L49").Value = "Especial"
$ws.Range("M49").Value = 330
$ws.Range("N49").Value = 8000
$ws.Range("O49").Value = 8000
$ws.Range("P49").Value = 8000
$ws.Range("R49").Value = "Provincia de Linares"
$ws.Range("S49").Value = 4000
$ws.Range("D50").Value2 = 44589
$ws.Range("M50").Value = 380
$ws.Range("N50").Value = 7000
$ws.Range("O50").Value = 7000
$ws.Range("P50").Value = 7000
$ws.Range("S50").Value = 3500
$ws.Range("D51").Value2 = 44238
$ws.Range("M51").Value = 200
$ws.Range("N51").Value = 6000
$ws.Range("O51").Value = 6000
$ws.Range("P51").Value = 6000
$ws.Range("S51").Value = 3000
$ws.Range("D52").Value2 = 44238
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 150
$ws.Range("N52").Value = 6000
$ws.Range("O52").Value = 6000
$ws.Range("P52").Value = 6000
$ws.Range("S52").Value = 3000
$ws.Range("D53").Value2 = 44356
$ws.Range("M53").Value = 60
$ws.Range("N53").Value = 10000
$ws.Range("O53").Value = 10000
$ws.Range("P53").Value = 10000
$ws.Range("R53").Value = "Provincia de Curicó"
$ws.Range("S53").Value = 5000
$ws.Range("D54").Value2 = 44582
$ws.Range("L54").Value = "Especial"
$ws.Range("M54").Value = 440
$ws.Range("N54").Value = 8000
$ws.Range("O54").Value = 8000
$ws.Range("P54").Value = 8000
$ws.Range("R54").Value = "Provincia de Linares"
$ws.Range("S54").Value = 4000
$ws.Range("D55").Value2 = 44582
$ws.Range("M55").Value = 380
$ws.Range("N55").Value = 7000
$ws.Range("O55").Value = 7000
$ws.Range("P55").Value = 7000
$ws.Range("R55").Value = "Provincia de Linares"
$ws.Range("S55").Value = 3500
$ws.Range("D56").Value2 = 44362
$ws.Range("M56").Value = 75
$ws.Range("N56").Value = 10000
$ws.Range("O56").Value = 10000
$ws.Range("P56").Value = 10000
$ws.Range("R56").Value = "Provincia de Curicó"
$ws.Range("S56").Value = 5000
$ws.Range("D57").Value2 = 44364
$ws.Range("M57").Value = 75
$ws.Range("N57").Value = 10000
$ws.Range("O57").Value = 10000
$ws.Range("P57").Value = 10000
$ws.Range("R57").Value = "Provincia de Curicó"
$ws.Range("S57").Value = 5000
$ws.Range("D58").Value2 = 44279
$ws.Range("M58").Value = 150
$ws.Range("N58").Value = 8000
$ws.Range("O58").Value = 8000
$ws.Range("P58").Value = 8000
$ws.Range("S58").Value = 4000
$ws.Range("D59").Value2 = 44552
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 610
$ws.Range("N59").Value = 7500
$ws.Range("P59").Value = 7730
$ws.Range("R59").Value = "Región del Maule"
$ws.Range("S59").Value = 3865
$ws.Range("D60").Value2 = 44217
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 250
$ws.Range("N60").Value = 6500
$ws.Range("O60").Value = 6600
$ws.Range("P60").Value = 6560
$ws.Range("R60").Value = "Provincia de Linares"
$ws.Range("S60").Value = 3280
$ws.Range("L61").Value = "Especial"
$ws.Range("M61").Value = 50
$ws.Range("N61").Value = 8000
$ws.Range("O61").Value = 8000
$ws.Range("P61").Value = 8000
$ws.Range("S61").Value = 4000
$ws.Range("D62").Value2 = 44166
$ws.Range("M62").Value = 150
$ws.Range("N62").Value = 7200
$ws.Range("O62").Value = 7200
$ws.Range("P62").Value = 7200
$ws.Range("R62").Value = "Región de O'Higgins"
$ws.Range("S62").Value = 3600
$ws.Range("D63").Value2 = 44166
$ws.Range("M63").Value = 80
$ws.Range("L64").Value = "Especial"
$ws.Range("M64").Value = 280
$ws.Range("N64").Value = 8000
$ws.Range("O64").Value = 8000
$ws.Range("P64").Value = 8000
$ws.Range("S64").Value = 4000
$ws.Range("D65").Value2 = 44581
$ws.Range("M65").Value = 250
$ws.Range("N65").Value = 7000
$ws.Range("O65").Value = 7000
$ws.Range("P65").Value = 7000
$ws.Range("S65").Value = 3500
$ws.Range("D66").Value2 = 44581
$ws.Range("L66").Value = "Segunda"
$ws.Range("M66").Value = 220
$ws.Range("N66").Value = 6000
$ws.Range("O66").Value = 6000
$ws.Range("P66").Value = 6000
$ws.Range("S66").Value = 3000
$ws.Range("D67").Value2 = 44573
$ws.Range("L67").Value = "Primera"
$ws.Range("M67").Value = 890
$ws.Range("N67").Value = 7500
$ws.Range("P67").Value = 7730
$ws.Range("S67").Value = 3865
$ws.Range("D68").Value2 = 44280
$ws.Range("M68").Value = 260
$ws.Range("N68").Value = 8000
$ws.Range("O68").Value = 8000
$ws.Range("P68").Value = 8000
$ws.Range("S68").Value = 4000
$ws.Range("D69").Value2 = 44588
$ws.Range("L69").Value = "Especial"
$ws.Range("M69").Value = 310
$ws.Range("N69").Value = 8000
$ws.Range("O69").Value = 8000
$ws.Range("P69").Value = 8000
$ws.Range("S69").Value = 4000
$ws.Range("D70").Value2 = 44588
$ws.Range("L70").Value = "Primera"
$ws.Range("M70").Value = 350
$ws.Range("N70").Value = 7000
$ws.Range("O70").Value = 7000
$ws.Range("P70").Value = 7000
$ws.Range("S70").Value = 3500
$ws.Range("D71").Value2 = 44187
$ws.Range("M71").Value = 220
$ws.Range("N71").Value = 7000
$ws.Range("P71").Value = 7000
$ws.Range("S71").Value = 3500
$ws.Range("D72").Value2 = 44187
$ws.Range("L72").Value = "Segunda"
$ws.Range("N72").Value = 5000
$ws.Range("O72").Value = 5000
$ws.Range("P72").Value = 5000
$ws.Range("S72").Value = 2500
$ws.Range("D73").Value2 = 44202
$ws.Range("M73").Value = 310
$ws.Range("N73").Value = 6500
$ws.Range("O73").Value = 7000
$ws.Range("P73").Value = 6677
$ws.Range("S73").Value = 3338
$ws.Range("D74").Value2 = 44225
$ws.Range("M74").Value = 260
$ws.Range("R74").Value = "Provincia de Linares"
$ws.Range("D75").Value2 = 44561
$ws.Range("M75").Value = 300
$ws.Range("N75").Value = 8000
$ws.Range("O75").Value = 8000
$ws.Range("P75").Value = 8000
$ws.Range("R75").Value = "Provincia de Linares"
$ws.Range("S75").Value = 4000
$ws.Range("D76").Value2 = 44264
$ws.Range("M76").Value = 220
$ws.Range("N76").Value = 6000
$ws.Range("O76").Value = 6000
$ws.Range("P76").Value = 6000
$ws.Range("R76").Value = "Región de O'Higgins"
$ws.Range("S76").Value = 3000
$ws.Range("D77").Value2 = 44533
$ws.Range("M77").Value = 350
$ws.Range("N77").Value = 10000
$ws.Range("O77").Value = 10000
$ws.Range("P77").Value = 10000
$ws.Range("S77").Value = 5000
$ws.Range("D78").Value2 = 44300
$ws.Range("D79").Value2 = 44306
$ws.Range("M79").Value = 200

# --- Step 2: insert two new rows before the old final row (80 -> 82) ---
$ws.Rows("80:81").Insert()

# --- Step 3: populate new row 80 ---
$ws.Range("A80").Value = 9
$ws.Range("B80").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C80").Value = "Metropolitana"
$ws.Range("D80").Value2 = 44299
$ws.Range("E80").Value = 13
$ws.Range("F80").Value = "Fruta"
$ws.Range("G80").Value = 100101
$ws.Range("H80").Value = "Berries"
$ws.Range("I80").Value = 100101004
$ws.Range("J80").Value = "Frambuesa"
$ws.Range("K80").Value = "Sin especificar"
$ws.Range("L80").Value = "Primera"
$ws.Range("M80").Value = 250
$ws.Range("N80").Value = 7000
$ws.Range("O80").Value = 7000
$ws.Range("P80").Value = 7000
$ws.Range("Q80").Value = "`$/bandeja 2 kilos"
$ws.Range("R80").Value = "Provincia de Curicó"
$ws.Range("S80").Value = 3500
$ws.Range("T80").Value = 2

# --- Step 4: populate new row 81 ---
$ws.Range("A81").Value = 9
$ws.Range("B81").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C81").Value = "Metropolitana"
$ws.Range("D81").Value2 = 44302
$ws.Range("E81").Value = 13
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100101
$ws.Range("H81").Value = "Berries"
$ws.Range("I81").Value = 100101004
$ws.Range("J81").Value = "Frambuesa"
$ws.Range("K81").Value = "Sin especificar"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 150
$ws.Range("N81").Value = 7000
$ws.Range("O81").Value = 7000
$ws.Range("P81").Value = 7000
$ws.Range("Q81").Value = "`$/bandeja 2 kilos"
$ws.Range("R81").Value = "Provincia de Curicó"
$ws.Range("S81").Value = 3500
$ws.Range("T81").Value = 2
